$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Role" column header and values to reflect the new
# PI / Sub I terminology (was Principal / Sub).
$ws.Range("A1").Value = "Role (PI/Sub I)"
$ws.Range("A2").Value = "PI"
$ws.Range("A3").Value = "Sub I"
$ws.Range("A4").Value = "Sub I"
$ws.Range("A5").Value = "Sub I"
